# "using gains for all"
# Inserts two new metric columns (M_TotalTax, M_CorpTax) right after M_POP,
# pushing the existing GFA/IMF/OECD columns two slots to the right, fills in
# the new columns' values for every group row, fixes the row-3 (LICs) cells
# that used to be the literal "inf" placeholder but now have real numbers in
# the two new columns, and corrects the UMICs M_POP value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before the old column F (GFA - Sales). Excel
# shifts the old F:M (GFA...OECD - Sales + Emp) block to H:O and carries the
# header/column formatting along with it.
$ws.Range("F1:G1").EntireColumn.Insert()

# New headers
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# High Income row
$ws.Range("F2").Value = 14106286460237.92
$ws.Range("G2").Value = 1155021202746.413

# LICs row - previously all of F3:M3 were the literal "inf"; now the two new
# columns carry real numbers while the (shifted) rest stay "inf".
$ws.Range("F3").Value = 3207987015.574299
$ws.Range("G3").Value = 0

# LMICs row
$ws.Range("F4").Value = 734615892234.8064
$ws.Range("G4").Value = 88889835996.30263

# Tax haven row
$ws.Range("F5").Value = 558865056646.082
$ws.Range("G5").Value = 72600947639.16805

# UMICs row
$ws.Range("F6").Value = 4579473077980.816
$ws.Range("G6").Value = 674619880691.7614

# UMICs M_POP value correction
$ws.Range("E6").Value = 2427884184.75
